$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16-18: keep worker PPT/862944/JOSELYN rows, shift "Periodo Mora" values
# and bump the "Salario Basico" (col G) from 1160000 to 1500000.
$ws.Range("E16").Value = "2308"
$ws.Range("G16").Value = 1500000

$ws.Range("E17").Value = "2309"
$ws.Range("G17").Value = 1500000

$ws.Range("E18").Value = "2310"
$ws.Range("G18").Value = 1500000

# Row 19 becomes the MELISA/CC record (moved up from row 22), value mora 2000.
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1007786694"
$ws.Range("D19").Value = "MELISA CHELIAN ALFARO PACHECO"
$ws.Range("E19").Value = "2310"
$ws.Range("F19").Value = 2000
$ws.Range("G19").Value = 1500000

# Rows 20-21: continue the PPT/JOSELYN series with later periods, G bumped.
$ws.Range("E20").Value = "2311"
$ws.Range("G20").Value = 1500000

$ws.Range("E21").Value = "2312"
$ws.Range("G21").Value = 1500000

# Row 22 becomes the PPT/JOSELYN record for period 2401 (moved down from row 16).
$ws.Range("B22").Value = "PPT"
$ws.Range("C22").Value = "862944"
$ws.Range("D22").Value = "JOSELYN ANTONIA FERNANDEZ PERNALETE"
$ws.Range("E22").Value = "2401"
$ws.Range("F22").Value = 46400
$ws.Range("G22").Value = 1500000
